$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.0296877408691
$ws.Range("D2").Value = 1.048562058525631
$ws.Range("E2").Value = 1.029488945739532
$ws.Range("F2").Value = 1.054722942665752
$ws.Range("I2").Value = 1.040465465200532
$ws.Range("J2").Value = 1.034833134630684
$ws.Range("K2").Value = 1.051321557922628
$ws.Range("L2").Value = 1.032302709353689
$ws.Range("M2").Value = 1.057465368533453
$ws.Range("N2").Value = 1.036302716378879
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.030542318384428
$ws.Range("D3").Value = 1.049148923898365
$ws.Range("E3").Value = 1.030212570071572
$ws.Range("F3").Value = 1.055466030677948
$ws.Range("I3").Value = 1.040644306336934
$ws.Range("J3").Value = 1.035329334738318
$ws.Range("K3").Value = 1.051721076403784
$ws.Range("L3").Value = 1.032834717319635
$ws.Range("M3").Value = 1.058021935620201
$ws.Range("N3").Value = 1.036799621147581
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.031095998679723
$ws.Range("D4").Value = 1.049529206289801
$ws.Range("E4").Value = 1.030681795685327
$ws.Range("F4").Value = 1.055947777522388
$ws.Range("I4").Value = 1.040759198702706
$ws.Range("J4").Value = 1.035650473993817
$ws.Range("K4").Value = 1.051979431694966
$ws.Range("L4").Value = 1.033179284395125
$ws.Range("M4").Value = 1.058382321403789
$ws.Range("N4").Value = 1.037121216457658
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.031328934614357
$ws.Range("D5").Value = 1.049689204654618
$ws.Range("E5").Value = 1.030879293884515
$ws.Range("F5").Value = 1.056150521844692
$ws.Range("I5").Value = 1.040807300037079
$ws.Range("J5").Value = 1.035785495128956
$ws.Range("K5").Value = 1.05208800459538
$ws.Range("L5").Value = 1.033324216189042
$ws.Range("M5").Value = 1.058533885654254
$ws.Range("N5").Value = 1.037256429338296
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.031368055431349
$ws.Range("D6").Value = 1.049716076513233
$ws.Range("E6").Value = 1.030912468473299
$ws.Range("F6").Value = 1.056184576227697
$ws.Range("I6").Value = 1.040815364746895
$ws.Range("J6").Value = 1.035808166566857
$ws.Range("K6").Value = 1.052106232096806
$ws.Range("L6").Value = 1.033348555279667
$ws.Range("M6").Value = 1.058559337314376
$ws.Range("N6").Value = 1.03727913297224
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.031099110520046
$ws.Range("D7").Value = 1.049531343696647
$ws.Range("E7").Value = 1.030684433741857
$ws.Range("F7").Value = 1.055950485748213
$ws.Range("I7").Value = 1.040759842219605
$ws.Range("J7").Value = 1.035652278098666
$ws.Range("K7").Value = 1.051980882608972
$ws.Range("L7").Value = 1.033181220684751
$ws.Range("M7").Value = 1.058384346386698
$ws.Range("N7").Value = 1.037123023124542
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.029976401060319
$ws.Range("D8").Value = 1.048760279196969
$ws.Range("E8").Value = 1.029733291427599
$ws.Range("F8").Value = 1.05497388114037
$ws.Range("I8").Value = 1.040526076760778
$ws.Range("J8").Value = 1.035000813636594
$ws.Range("K8").Value = 1.051456609334037
$ws.Range("L8").Value = 1.032482436365715
$ws.Range("M8").Value = 1.057653410268836
$ws.Range("N8").Value = 1.036470633508211
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.028003563477567
$ws.Range("D9").Value = 1.047405795643958
$ws.Range("E9").Value = 1.028064936607431
$ws.Range("F9").Value = 1.053260112504269
$ws.Range("I9").Value = 1.040107830418167
$ws.Range("J9").Value = 1.033853401221567
$ws.Range("K9").Value = 1.050531609720553
$ws.Range("L9").Value = 1.03125361873439
$ws.Range("M9").Value = 1.056367394839282
$ws.Range("N9").Value = 1.035321591635944
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.026692144849693
$ws.Range("D10").Value = 1.046505766371479
$ws.Range("E10").Value = 1.026957971013343
$ws.Range("F10").Value = 1.052122522513225
$ws.Range("I10").Value = 1.039824794554424
$ws.Range("J10").Value = 1.033088906278793
$ws.Range("K10").Value = 1.049914244285119
$ws.Range("L10").Value = 1.030436189844783
$ws.Range("M10").Value = 1.055511489210072
$ws.Range("N10").Value = 1.034556011022663
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.026125208941538
$ws.Range("D11").Value = 1.046116770238265
$ws.Range("E11").Value = 1.026479915679503
$ws.Range("F11").Value = 1.051631127794532
$ws.Range("I11").Value = 1.039701249134659
$ws.Range("J11").Value = 1.032757993522537
$ws.Range("K11").Value = 1.049646770012971
$ws.Range("L11").Value = 1.03008267432479
$ws.Range("M11").Value = 1.055141235251163
$ws.Range("N11").Value = 1.034224628332337
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.025914763076489
$ws.Range("D12").Value = 1.045972390292208
$ws.Range("E12").Value = 1.026302537009608
$ws.Range("F12").Value = 1.051448782769731
$ws.Range("I12").Value = 1.039655211038931
$ws.Range("J12").Value = 1.032635096900413
$ws.Range("K12").Value = 1.049547396900192
$ws.Range("L12").Value = 1.029951429913476
$ws.Range("M12").Value = 1.055003762103234
$ws.Range("N12").Value = 1.034101557182912
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.025959898096159
$ws.Range("D13").Value = 1.046003355246716
$ws.Range("E13").Value = 1.026340576601704
$ws.Range("F13").Value = 1.051487888177564
$ws.Range("I13").Value = 1.039665093043906
$ws.Range("J13").Value = 1.032661457757559
$ws.Range("K13").Value = 1.049568713709404
$ws.Range("L13").Value = 1.029979579231147
$ws.Range("M13").Value = 1.055033248013434
$ws.Range("N13").Value = 1.034127955475499
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.026107810568206
$ws.Range("D14").Value = 1.046104833480433
$ws.Range("E14").Value = 1.026465249568387
$ws.Range("F14").Value = 1.051616051395947
$ws.Range("I14").Value = 1.039697446622487
$ws.Range("J14").Value = 1.032747834452117
$ws.Range("K14").Value = 1.049638556225924
$ws.Range("L14").Value = 1.030071824246066
$ws.Range("M14").Value = 1.055129870532359
$ws.Range("N14").Value = 1.034214454834872
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.026198962840418
$ws.Range("D15").Value = 1.046167372276319
$ws.Range("E15").Value = 1.026542090243713
$ws.Range("F15").Value = 1.051695041014691
$ws.Range("I15").Value = 1.039717361164545
$ws.Range("J15").Value = 1.032801056555249
$ws.Range("K15").Value = 1.049681585730954
$ws.Range("L15").Value = 1.030128668363376
$ws.Range("M15").Value = 1.055189410282702
$ws.Range("N15").Value = 1.034267752519495
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.026729789975906
$ws.Range("D16").Value = 1.046531598155801
$ws.Range("E16").Value = 1.026989724846449
$ws.Range("F16").Value = 1.052155160004795
$ws.Range("I16").Value = 1.039832973084612
$ws.Range("J16").Value = 1.033110870477541
$ws.Range("K16").Value = 1.049931992612412
$ws.Range("L16").Value = 1.030459660823506
$ws.Range("M16").Value = 1.055536069455262
$ws.Range("N16").Value = 1.034578006413092
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.027063010454328
$ws.Range("D17").Value = 1.046760262264402
$ws.Range("E17").Value = 1.027270855014045
$ws.Range("F17").Value = 1.05244410042879
$ws.Range("I17").Value = 1.039905229136641
$ws.Range("J17").Value = 1.033305241229458
$ws.Range("K17").Value = 1.050089026820187
$ws.Range("L17").Value = 1.030667401524158
$ws.Range("M17").Value = 1.05575361689838
$ws.Range("N17").Value = 1.034772653193772
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.027257460692316
$ws.Range("D18").Value = 1.046893707745757
$ws.Range("E18").Value = 1.027434955710961
$ws.Range("F18").Value = 1.052612749089301
$ws.Range("I18").Value = 1.039947279391645
$ws.Range("J18").Value = 1.033418625756992
$ws.Range("K18").Value = 1.050180607555594
$ws.Range("L18").Value = 1.030788615190153
$ws.Range("M18").Value = 1.055880543145745
$ws.Range("N18").Value = 1.03488619874034
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.027323778138892
$ws.Range("D19").Value = 1.046939220946868
$ws.Range("E19").Value = 1.027490930483213
$ws.Range("F19").Value = 1.052670273295964
$ws.Range("I19").Value = 1.039961601224049
$ws.Range("J19").Value = 1.033457288837736
$ws.Range("K19").Value = 1.05021183168964
$ws.Range("L19").Value = 1.030829953004196
$ws.Range("M19").Value = 1.055923827515362
$ws.Range("N19").Value = 1.034924916727094
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.027027249898408
$ws.Range("D20").Value = 1.046735721558709
$ws.Range("E20").Value = 1.027240679751825
$ws.Range("F20").Value = 1.052413087986489
$ws.Range("I20").Value = 1.039897486616283
$ws.Range("J20").Value = 1.033284385903319
$ws.Range("K20").Value = 1.05007218004207
$ws.Range("L20").Value = 1.030645108555546
$ws.Range("M20").Value = 1.05573027253502
$ws.Range("N20").Value = 1.034751768250677
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.026064250171506
$ws.Range("D21").Value = 1.046074947602511
$ws.Range("E21").Value = 1.026428531161963
$ws.Range("F21").Value = 1.051578305512587
$ws.Range("I21").Value = 1.039687923374568
$ws.Range("J21").Value = 1.032722398130566
$ws.Range("K21").Value = 1.049617989921024
$ws.Range("L21").Value = 1.03004465852773
$ws.Range("M21").Value = 1.055101416058326
$ws.Range("N21").Value = 1.034188982390826
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.025459581194155
$ws.Range("D22").Value = 1.045660133311773
$ws.Range("E22").Value = 1.025919014998799
$ws.Range("F22").Value = 1.05105449207652
$ws.Range("I22").Value = 1.039555307750547
$ws.Range("J22").Value = 1.032369165313628
$ws.Range("K22").Value = 1.049332300271029
$ws.Range("L22").Value = 1.029667520058322
$ws.Range("M22").Value = 1.054706351916227
$ws.Range("N22").Value = 1.033835247942772
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.025780050612574
$ws.Range("D23").Value = 1.045879972837344
$ws.Range("E23").Value = 1.026189012955112
$ws.Range("F23").Value = 1.051332075477311
$ws.Range("I23").Value = 1.039625690580956
$ws.Range("J23").Value = 1.032556409741335
$ws.Range("K23").Value = 1.049483760877572
$ws.Range("L23").Value = 1.029867411022906
$ws.Range("M23").Value = 1.054915751703785
$ws.Range("N23").Value = 1.034022758279044
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.027043408281474
$ws.Range("D24").Value = 1.046746810232561
$ws.Range("E24").Value = 1.027254314276526
$ws.Range("F24").Value = 1.052427100820448
$ws.Range("I24").Value = 1.039900985423067
$ws.Range("J24").Value = 1.033293809492566
$ws.Range("K24").Value = 1.050079792421428
$ws.Range("L24").Value = 1.030655181658785
$ws.Range("M24").Value = 1.055740820741593
$ws.Range("N24").Value = 1.034761205222502
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.028512925354234
$ws.Range("D25").Value = 1.04775544886148
$ws.Range("E25").Value = 1.028495325159521
$ws.Range("F25").Value = 1.053702304401517
$ws.Range("I25").Value = 1.040216701608164
$ws.Range("J25").Value = 1.034149961678412
$ws.Range("K25").Value = 1.050770873572829
$ws.Range("L25").Value = 1.031570989314725
$ws.Range("M25").Value = 1.05669961464451
$ws.Range("N25").Value = 1.035618573242654
